$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Total Holding"
$ws.Range("D1").Value = "Total Balance"
